# The sheet lists one "label_*.fasta" record per row (column A), with constant
# columns B/C/D. The commit shifts every data row from 66..323 up by 64 rows
# (i.e. row r's new value is the old value that used to live at row r+64),
# then drops the now-surplus trailing rows 324..387, shrinking the used range
# from A1:D387 to A1:D323. Rows 2..65 are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$shift = 64
$startRow = 66
$endRow = 323
$lastRow = 387

# Read every source value (row + 64) up front, before any cell in that same
# range gets overwritten, since the read/write windows overlap.
$vals = @{}
for ($r = $startRow; $r -le $endRow; $r++) {
    $vals[$r] = $ws.Cells.Item($r + $shift, 1).Value2
}

for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $vals[$r]
}

# Remove the trailing rows that are no longer needed (324..387), shrinking
# the sheet's dimension to A1:D323.
$ws.Range("A" + ($endRow + 1) + ":D" + $lastRow).EntireRow.Delete()
